$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*OpenNeuro PET, Stanford University*") {
        $p.Range.Delete()
        break
    }
}
